$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data: Player, Position, Team for rows 2-19
$data = @(
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Andrew Wiggins", "SF,PF", "Miami Heat"),
    @("Brandon Clarke", "PF,C", "Memphis Grizzlies"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Zach Edey", "C", "Memphis Grizzlies"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Jordan Hawkins", "SG", "New Orleans Pelicans"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Ausar Thompson", "SF,PF", "Detroit Pistons"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Anthony Davis", "PF,C", "Dallas Mavericks"),
    @("Brandon Williams", "PG", "Dallas Mavericks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
